$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update swapped / rotated match rows (only cells whose value actually changes) ---
$ws.Range("B27").Value = 6504313
$ws.Range("F27").Value = "Guabira"
$ws.Range("G27").Value = "Atletico Palmaflor Vinto"
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = "H"
$ws.Range("K27").Value = 1.75
$ws.Range("M27").Value = 4
$ws.Range("N27").Value = 2
$ws.Range("O27").Value = 3.6
$ws.Range("P27").Value = 3.6
$ws.Range("Q27").Value = -0.25
$ws.Range("R27").Value = 1.75
$ws.Range("S27").Value = 2.05
$ws.Range("T27").Value = 2.5
$ws.Range("U27").Value = 1.85
$ws.Range("V27").Value = 1.95
$ws.Range("W27").Value = 1
$ws.Range("Y27").Value = -1
$ws.Range("Z27").Value = 0.75
$ws.Range("AA27").Value = -1
$ws.Range("AB27").Value = -1
$ws.Range("AC27").Value = 0.95
$ws.Range("B28").Value = 6504831
$ws.Range("F28").Value = "Libertad Gran Mamore FC"
$ws.Range("G28").Value = "Always Ready"
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = "A"
$ws.Range("K28").Value = 3.3
$ws.Range("M28").Value = 1.909
$ws.Range("N28").Value = 3.1
$ws.Range("O28").Value = 3.5
$ws.Range("P28").Value = 2.2
$ws.Range("Q28").Value = 0.25
$ws.Range("R28").Value = 1.9
$ws.Range("S28").Value = 1.9
$ws.Range("T28").Value = 2.75
$ws.Range("U28").Value = 2.025
$ws.Range("V28").Value = 1.775
$ws.Range("W28").Value = -1
$ws.Range("Y28").Value = 1.2
$ws.Range("Z28").Value = -1
$ws.Range("AA28").Value = 0.8999999999999999
$ws.Range("AB28").Value = 0.5125
$ws.Range("AC28").Value = -0.5
$ws.Range("B47").Value = 6504291
$ws.Range("F47").Value = "Club Aurora"
$ws.Range("G47").Value = "Universitario De Vinto"
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 1
$ws.Range("M47").Value = 3.3
$ws.Range("N47").Value = 1.7
$ws.Range("P47").Value = 5.5
$ws.Range("R47").Value = 1.85
$ws.Range("S47").Value = 1.95
$ws.Range("U47").Value = 1.9
$ws.Range("V47").Value = 1.9
$ws.Range("Y47").Value = 4.5
$ws.Range("AA47").Value = 0.95
$ws.Range("AB47").Value = -1
$ws.Range("AC47").Value = 0.8999999999999999
$ws.Range("B48").Value = 6504578
$ws.Range("F48").Value = "Guabira"
$ws.Range("G48").Value = "Real Santa Cruz"
$ws.Range("H48").Value = 1
$ws.Range("I48").Value = 2
$ws.Range("M48").Value = 3.4
$ws.Range("N48").Value = 1.75
$ws.Range("P48").Value = 5.25
$ws.Range("R48").Value = 1.95
$ws.Range("S48").Value = 1.85
$ws.Range("U48").Value = 2
$ws.Range("V48").Value = 1.8
$ws.Range("Y48").Value = 4.25
$ws.Range("AA48").Value = 0.8500000000000001
$ws.Range("AB48").Value = 1
$ws.Range("AC48").Value = -1
$ws.Range("B128").Value = 7462738
$ws.Range("F128").Value = "Vaca Diez"
$ws.Range("G128").Value = "The Strongest"
$ws.Range("H128").Value = 2
$ws.Range("I128").Value = 2
$ws.Range("J128").Value = "D"
$ws.Range("K128").Value = 4
$ws.Range("L128").Value = 4
$ws.Range("M128").Value = 1.666
$ws.Range("N128").Value = 4
$ws.Range("O128").Value = 3.8
$ws.Range("P128").Value = 1.75
$ws.Range("Q128").Value = 0.75
$ws.Range("R128").Value = 1.8
$ws.Range("S128").Value = 2
$ws.Range("T128").Value = 3
$ws.Range("U128").Value = 1.925
$ws.Range("V128").Value = 1.875
$ws.Range("W128").Value = -1
$ws.Range("X128").Value = 2.8
$ws.Range("Z128").Value = 0.8
$ws.Range("AB128").Value = 0.925
$ws.Range("AC128").Value = -1
$ws.Range("B129").Value = 7462542
$ws.Range("F129").Value = "Always Ready"
$ws.Range("G129").Value = "Royal Pari FC"
$ws.Range("H129").Value = 3
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = "H"
$ws.Range("K129").Value = 1.363
$ws.Range("L129").Value = 4.75
$ws.Range("M129").Value = 6.5
$ws.Range("N129").Value = 1.285
$ws.Range("O129").Value = 6.5
$ws.Range("P129").Value = 8
$ws.Range("Q129").Value = -1.75
$ws.Range("R129").Value = 1.9
$ws.Range("S129").Value = 1.9
$ws.Range("T129").Value = 3.25
$ws.Range("U129").Value = 1.85
$ws.Range("V129").Value = 1.95
$ws.Range("W129").Value = 0.2849999999999999
$ws.Range("X129").Value = -1
$ws.Range("Z129").Value = 0.8999999999999999
$ws.Range("AB129").Value = -0.5
$ws.Range("AC129").Value = 0.475
$ws.Range("B149").Value = 7532420
$ws.Range("F149").Value = "Club Aurora"
$ws.Range("G149").Value = "Vaca Diez"
$ws.Range("H149").Value = 3
$ws.Range("K149").Value = 1.333
$ws.Range("L149").Value = 5
$ws.Range("M149").Value = 8
$ws.Range("N149").Value = 1.3
$ws.Range("O149").Value = 6.5
$ws.Range("P149").Value = 7
$ws.Range("R149").Value = 1.8
$ws.Range("S149").Value = 2
$ws.Range("T149").Value = 3.25
$ws.Range("U149").Value = 1.95
$ws.Range("V149").Value = 1.85
$ws.Range("W149").Value = 0.3
$ws.Range("Z149").Value = 0.8
$ws.Range("AB149").Value = -0.5
$ws.Range("AC149").Value = 0.425
$ws.Range("B150").Value = 7532421
$ws.Range("F150").Value = "Guabira"
$ws.Range("G150").Value = "Independiente Petrolero"
$ws.Range("H150").Value = 2
$ws.Range("K150").Value = 1.4
$ws.Range("L150").Value = 4.5
$ws.Range("M150").Value = 7.5
$ws.Range("N150").Value = 1.333
$ws.Range("O150").Value = 5.5
$ws.Range("P150").Value = 9.5
$ws.Range("R150").Value = 1.85
$ws.Range("S150").Value = 1.95
$ws.Range("T150").Value = 3
$ws.Range("U150").Value = 1.825
$ws.Range("V150").Value = 1.975
$ws.Range("W150").Value = 0.333
$ws.Range("Z150").Value = 0.8500000000000001
$ws.Range("AB150").Value = -1
$ws.Range("AC150").Value = 0.9750000000000001
$ws.Range("B153").Value = 7532431
$ws.Range("F153").Value = "Blooming"
$ws.Range("G153").Value = "Always Ready"
$ws.Range("H153").Value = 2
$ws.Range("I153").Value = 1
$ws.Range("K153").Value = 2.3
$ws.Range("L153").Value = 3.5
$ws.Range("M153").Value = 2.625
$ws.Range("N153").Value = 1.833
$ws.Range("O153").Value = 4
$ws.Range("P153").Value = 3.6
$ws.Range("Q153").Value = -0.5
$ws.Range("R153").Value = 1.825
$ws.Range("S153").Value = 1.975
$ws.Range("T153").Value = 3
$ws.Range("U153").Value = 2
$ws.Range("V153").Value = 1.8
$ws.Range("W153").Value = 0.833
$ws.Range("Z153").Value = 0.825
$ws.Range("AB153").Value = 0
$ws.Range("AC153").Value = -0
$ws.Range("B154").Value = 7532417
$ws.Range("F154").Value = "Real Tomayapo"
$ws.Range("G154").Value = "Atletico Palmaflor Vinto"
$ws.Range("H154").Value = 4
$ws.Range("I154").Value = 0
$ws.Range("K154").Value = 1.3
$ws.Range("L154").Value = 4.5
$ws.Range("M154").Value = 8.5
$ws.Range("N154").Value = 1.166
$ws.Range("O154").Value = 8
$ws.Range("P154").Value = 12
$ws.Range("Q154").Value = -2.25
$ws.Range("R154").Value = 1.95
$ws.Range("S154").Value = 1.85
$ws.Range("T154").Value = 3.75
$ws.Range("U154").Value = 1.975
$ws.Range("V154").Value = 1.825
$ws.Range("W154").Value = 0.1659999999999999
$ws.Range("Z154").Value = 0.95
$ws.Range("AB154").Value = 0.4875
$ws.Range("AC154").Value = -0.5
$ws.Range("B142").Value = 7532430
$ws.Range("F142").Value = "Always Ready"
$ws.Range("G142").Value = "Oriente Petrolero"
$ws.Range("H142").Value = 4
$ws.Range("I142").Value = 1
$ws.Range("K142").Value = 1.4
$ws.Range("L142").Value = 4.2
$ws.Range("M142").Value = 7
$ws.Range("N142").Value = 1.363
$ws.Range("O142").Value = 4.5
$ws.Range("P142").Value = 8.5
$ws.Range("Q142").Value = -1.5
$ws.Range("U142").Value = 1.9
$ws.Range("V142").Value = 1.9
$ws.Range("W142").Value = 0.363
$ws.Range("Z142").Value = 1
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.8999999999999999
$ws.Range("AC142").Value = -1
$ws.Range("B143").Value = 7532414
$ws.Range("F143").Value = "Independiente Petrolero"
$ws.Range("G143").Value = "Real Santa Cruz"
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 0
$ws.Range("J143").Value = "H"
$ws.Range("K143").Value = 1.571
$ws.Range("L143").Value = 3.75
$ws.Range("M143").Value = 5
$ws.Range("N143").Value = 1.3
$ws.Range("O143").Value = 5
$ws.Range("P143").Value = 11
$ws.Range("Q143").Value = -1.75
$ws.Range("R143").Value = 2
$ws.Range("S143").Value = 1.8
$ws.Range("T143").Value = 3
$ws.Range("U143").Value = 1.85
$ws.Range("V143").Value = 1.95
$ws.Range("W143").Value = 0.3
$ws.Range("Y143").Value = -1
$ws.Range("AA143").Value = 0.8
$ws.Range("AC143").Value = 0.95
$ws.Range("B144").Value = 7532413
$ws.Range("F144").Value = "Libertad Gran Mamore FC"
$ws.Range("G144").Value = "Club Aurora"
$ws.Range("H144").Value = 0
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 2.25
$ws.Range("L144").Value = 3.3
$ws.Range("M144").Value = 2.8
$ws.Range("N144").Value = 2.375
$ws.Range("O144").Value = 3.4
$ws.Range("P144").Value = 2.875
$ws.Range("Q144").Value = -0.25
$ws.Range("R144").Value = 2.025
$ws.Range("S144").Value = 1.775
$ws.Range("T144").Value = 2.5
$ws.Range("W144").Value = -1
$ws.Range("Y144").Value = 1.875
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.7749999999999999
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.8999999999999999

# --- Remove the two trailing (not-yet-played) fixtures, rows 219 and 220 ---
$ws.Range("A219:A220").EntireRow.Delete() | Out-Null

Write-Host "edit applied"
